$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.543.26"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.623.48"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'211.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "

$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'23.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("E9").Value = "  +1.52%  "

$ws.Range("D10").Value = "'0.0611"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.13%  "

$ws.Range("D11").Value = "'0.0878"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").Value = "1.853.13"
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("D13").Value = "1.634.71"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").Value = "'0.549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "

$ws.Range("D16").Value = "'65.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").Value = "27.499.15"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").Value = "'229.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").Value = "0.0₃0717"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "'10.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.54%  "

$ws.Range("D23").Value = "'4.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("E24").Value = "  +5.28%  "

$ws.Range("D25").Value = "'148.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("D26").Value = "'6.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.20%  "

$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "'15.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.89%  "

$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("D31").Value = "'0.0483"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "

$ws.Range("D32").Value = "'3.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").Value = "1.466.45"
$ws.Range("E33").Value = "  +1.45%  "

$ws.Range("E34").Value = "  -2.72%  "

$ws.Range("E35").Value = "  -2.69%  "

$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("D37").Value = "'0.939"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.74%  "

$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("D40").Value = "'0.553"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.80%  "

$ws.Range("E42").Value = "  -2.14%  "

$ws.Range("D43").Value = "'67.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.93%  "

$ws.Range("E44").Value = "  -0.76%  "

$ws.Range("E45").Value = "  -1.92%  "

$ws.Range("E46").Value = "  -6.03%  "

$ws.Range("D47").Value = "1.763.43"
$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("D49").Value = "'87.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.82%  "

$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").Value = "'0.0997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.81%  "
